# Update "want to go" counts (column F) per commit: regenerated gh-pages output at 456a3b4
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 2357   # was 2354
$ws1.Range("F8").Value = 1754   # was 1750
$ws1.Range("F9").Value = 2971   # was 2965
$ws1.Range("F11").Value = 4431   # was 4424
$ws1.Range("F12").Value = 388   # was 387
$ws1.Range("F13").Value = 213   # was 212
$ws1.Range("F15").Value = 558   # was 556
$ws1.Range("F18").Value = 218   # was 209
$ws1.Range("F20").Value = 108   # was 106
$ws1.Range("F21").Value = 306   # was 305
$ws1.Range("F22").Value = 4482   # was 4471
$ws1.Range("F23").Value = 7   # was 6
$ws1.Range("F24").Value = 3673   # was 3672
$ws1.Range("F25").Value = 1139   # was 1138
$ws1.Range("F26").Value = 215   # was 213
$ws1.Range("F27").Value = 567   # was 564
$ws1.Range("F30").Value = 586   # was 579
$ws1.Range("F31").Value = 580   # was 575
$ws1.Range("F32").Value = 543   # was 539

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F9").Value = 2357   # was 2354
$ws4.Range("F11").Value = 1754   # was 1750
$ws4.Range("F13").Value = 2971   # was 2965
$ws4.Range("F15").Value = 4431   # was 4424
$ws4.Range("F16").Value = 388   # was 387
$ws4.Range("F17").Value = 213   # was 212
$ws4.Range("F19").Value = 558   # was 556
$ws4.Range("F22").Value = 218   # was 209
$ws4.Range("F25").Value = 108   # was 106
$ws4.Range("F26").Value = 306   # was 305
$ws4.Range("F27").Value = 4482   # was 4471
$ws4.Range("F28").Value = 7   # was 6
$ws4.Range("F29").Value = 3673   # was 3672
$ws4.Range("F30").Value = 1139   # was 1138
$ws4.Range("F31").Value = 215   # was 213
$ws4.Range("F32").Value = 567   # was 564
$ws4.Range("F35").Value = 586   # was 579
$ws4.Range("F36").Value = 580   # was 575
$ws4.Range("F37").Value = 543   # was 539
